# Generate Report for Handoff
#
# The localization-status report is regenerated: the file
# "ab3a429d-9932-40a9-8d2d-981e3e741847.md" (still "In Translation") now
# appears first (row 2) on every sheet, and
# "178fcd04-bf82-41e8-b04a-dcfbf8504dc5.md" -- which was just produced for
# handoff -- now appears second (row 3) with a refreshed "Ready for
# handoff"/"In Translation"-follow-up status and new handoff timestamps.
#
# Hyperlink targets (the rIds / underlying URLs) are untouched; only the
# visible display text of each hyperlink is updated to track the new
# row/file pairing, exactly like the cell text it decorates.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "ab3a429d-9932-40a9-8d2d-981e3e741847.md"
$ws.Range("B2").Value = "In Translation"
$ws.Range("C2").Value = "In Translation"
$ws.Range("D2").Value = "2016-14-12 20:14:23"

$ws.Range("A3").Value = "178fcd04-bf82-41e8-b04a-dcfbf8504dc5.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "2016-15-12 20:15:43"

foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "ab3a429d-9932-40a9-8d2d-981e3e741847.md"
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = "178fcd04-bf82-41e8-b04a-dcfbf8504dc5.md"
    }
}

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "ab3a429d-9932-40a9-8d2d-981e3e741847.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "In Translation"
$ws.Range("D2").Value = "ab3a429d-9932-40a9-8d2d-981e3e741847.63e659546980afac8059c1fccd0a481757577c1e.zh-cn.xlf"
$ws.Range("E2").Value = "2016-03-12 20:13:35"
$ws.Range("H2").Value = "0001-01-01 00:00:00"
$ws.Range("I2").Value = "Include"

$ws.Range("A3").Value = "178fcd04-bf82-41e8-b04a-dcfbf8504dc5.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "178fcd04-bf82-41e8-b04a-dcfbf8504dc5.9dad81d3ebc399e4700936b0137ee8d9bbd47c0b.zh-cn.xlf"
$ws.Range("E3").Value = "2016-03-12 20:15:40"
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("I3").Value = "Include"

foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "ab3a429d-9932-40a9-8d2d-981e3e741847.md"
    } elseif ($addr -eq '$D$2') {
        $hl.TextToDisplay = "ab3a429d-9932-40a9-8d2d-981e3e741847.63e659546980afac8059c1fccd0a481757577c1e.zh-cn.xlf"
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = "178fcd04-bf82-41e8-b04a-dcfbf8504dc5.md"
    } elseif ($addr -eq '$D$3') {
        $hl.TextToDisplay = "178fcd04-bf82-41e8-b04a-dcfbf8504dc5.9dad81d3ebc399e4700936b0137ee8d9bbd47c0b.zh-cn.xlf"
    }
}

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "ab3a429d-9932-40a9-8d2d-981e3e741847.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "In Translation"
$ws.Range("D2").Value = "ab3a429d-9932-40a9-8d2d-981e3e741847.63e659546980afac8059c1fccd0a481757577c1e.de-de.xlf"
$ws.Range("E2").Value = "2016-03-12 20:14:23"
$ws.Range("H2").Value = "0001-01-01 00:00:00"
$ws.Range("I2").Value = "Include"

$ws.Range("A3").Value = "178fcd04-bf82-41e8-b04a-dcfbf8504dc5.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "178fcd04-bf82-41e8-b04a-dcfbf8504dc5.9dad81d3ebc399e4700936b0137ee8d9bbd47c0b.de-de.xlf"
$ws.Range("E3").Value = "2016-03-12 20:15:43"
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("I3").Value = "Include"

foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "ab3a429d-9932-40a9-8d2d-981e3e741847.md"
    } elseif ($addr -eq '$D$2') {
        $hl.TextToDisplay = "ab3a429d-9932-40a9-8d2d-981e3e741847.63e659546980afac8059c1fccd0a481757577c1e.de-de.xlf"
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = "178fcd04-bf82-41e8-b04a-dcfbf8504dc5.md"
    } elseif ($addr -eq '$D$3') {
        $hl.TextToDisplay = "178fcd04-bf82-41e8-b04a-dcfbf8504dc5.9dad81d3ebc399e4700936b0137ee8d9bbd47c0b.de-de.xlf"
    }
}
